$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "Yellowfin tuna "
$ws.Range("B7").Value = "Bluefin tuna"
$ws.Range("B10").Value = "Spiny lobster"
$ws.Range("B11").Value = "Abalone"
$ws.Range("B12").Value = "Bonito"
$ws.Range("B15").Value = "Rockfish"
$ws.Range("B22").Value = "California barracuda"
$ws.Range("B26").Value = "Spiny lobster"

$ws.Range("B30").Select()
